$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You are a conversion optimization analyst. You analyze the user recordings for users who visited the website, spent more than two minutes on the landing page, and did not submit their details in lead forms, which is your landing page goal. Observations:Users frequently use the chat functionality, which has a long wait time.The lead form is present above the first fold and has no button clicks.80% of users drop after scrolling 50% of the page.User reviews are at the bottom of the page.What changes should you make to get qualified leads?",
        "ques_type": 2,
        "options": [
            "Include frequently asked questions based on analysis of chat transcripts.",
            "Include a user form to initiate chat.",
            "Bring user reviews up to the first 50% of the page.",
            "Move the lead form to the bottom of the page."
        ],
        "score": "Include frequently asked questions based on analysis of chat transcripts."
    },
    {
        "title": "You are an e-commerce product strategist working with a fashion e-commerce store. While analyzing website data, you observe that the website product page has a 0% bounce rate and a 96% exit rate.  What does this say about the website\u2019s users?",
        "ques_type": 2,
        "options": [
            "Users have low intent and high willingness to buy.",
            "Users have low intent and low willingness to buy.",
            "Users have high intent and low willingness to buy.",
            "Users have high intent and high willingness to buy."
        ],
        "score": "Users have high intent and low willingness to buy."
    },
    {
        "title": "You are a landing page (LP) optimization specialist, and you are in the process of conducting a copy teardown for a company\u2019s LP. There is huge traffic coming in from a high intent Google search however, 97% of that traffic is bouncing at the first fold of the page, potentially due to a misleading headline. Which three key elements should you use to make headline copy changes if you\u2019re focusing on persuasive messaging and conversion?",
        "ques_type": 2,
        "options": [
            "Brand name, product/service price, customer pain points",
            "Authority, fear of missing out (FOMO), services",
            "Brand slogan, discounts, benefits",
            "Motivation, value proposition, incentive"
        ],
        "score": "Motivation, value proposition, incentive"
    },
    {
        "title": "As a conversion rate optimization analyst, you are analyzing a website landing page which takes inquiries for Web Development courses. You have found that users have a high number of fears, uncertainties, and doubts. You wish to address these by including important page elements like trust badges, guarantees, and other assurances in the highly engaging section of the page.  What type of user research tool should you use to find the highly engaging section of the page?",
        "ques_type": 2,
        "options": [
            "Google Analytics",
            "One-page user surveys",
            "Website heatmaps",
            "Scrollmap"
        ],
        "score": "Website heatmaps"
    }
]
'@

# The source cell (A1) previously held a placeholder 0 with a bold/bordered,
# centered style; the new layout puts the (now pretty-printed) questions
# string directly into A1 and drops the old A2 duplicate + the one-off style.
$ws.Range("A1").ClearFormats()
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $questionsText

# Re-fit the row so the tall placeholder-driven row height doesn't linger.
$ws.Rows(1).AutoFit()
